# The post previously on row 775 ("「休みが終わった」...") was removed from
# the source data. Delete that entire row so every subsequent post shifts
# up by one row (old row 776 -> 775, ..., old row 859 -> 858), matching the
# updated dimension ref of A1:C858.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(775).Delete()
